# HighLevelSequenceDiagrams.pptx edit:
#  - refresh the "datetimeFigureOut" date placeholders (slide master, all
#    slide layouts, notes master) from 9/21/2018 -> 11/10/2018
#  - update a handful of sequence-diagram text boxes on slide 1 to match
#    the current implementation (delete_friend -> delete, deletePerson ->
#    delete, EriumChangedEvent -> AddressBookChangedEvent)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders: slide master + every slide layout + notes master
# ---------------------------------------------------------------------

function Get-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14) {
            $phType = -1
            try { $phType = $shp.PlaceholderFormat.Type } catch {}
            if ($phType -eq 16) {
                return $shp
            }
        }
    }
    return $null
}

$newDate = "11/10/2018"

$sm = $p.SlideMaster

$masterDateShape = Get-DatePlaceholder $sm.Shapes
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $layout = $sm.CustomLayouts.Item($li)
    $layoutDateShape = Get-DatePlaceholder $layout.Shapes
    if ($layoutDateShape -ne $null) {
        $layoutDateShape.TextFrame.TextRange.Text = $newDate
    }
}

$nm = $p.NotesMaster
$notesDateShape = Get-DatePlaceholder $nm.Shapes
if ($notesDateShape -ne $null) {
    $notesDateShape.TextFrame.TextRange.Text = $newDate
}

# ---------------------------------------------------------------------
# 2) Slide 1 text updates
# ---------------------------------------------------------------------

$s = $p.Slides.Item(1)

# execute("delete_friend 1") -> execute("delete 1")
# Collapse the 3 runs into the first run's formatting.
$shp = $s.Shapes.Item("TextBox 25")
$tr = $shp.TextFrame.TextRange
$tail = $tr.Characters(10, $tr.Text.Length - 9)
$tail.Text = ""
$tr.Text = "execute(" + [char]8220 + "delete 1" + [char]8221 + ")"

# deletePerson(p) -> delete(p)
# Drop the first run entirely, keep the second run's formatting.
$shp = $s.Shapes.Item("TextBox 28")
$tr = $shp.TextFrame.TextRange
$head = $tr.Characters(1, 12)
$head.Text = ""
$tr.Text = "delete(p)"

# post(EriumChangedEvent) -> post(AddressBookChangedEvent)
$shp = $s.Shapes.Item("TextBox 32")
$tr = $shp.TextFrame.TextRange
$mid = $tr.Characters(6, 17)
$mid.Text = "AddressBookChangedEvent"

$shp = $s.Shapes.Item("TextBox 61")
$tr = $shp.TextFrame.TextRange
$mid = $tr.Characters(6, 17)
$mid.Text = "AddressBookChangedEvent"

# handleEriumChangedEvent() -> handleAddressBookChangedEvent()
$shp = $s.Shapes.Item("TextBox 73")
$tr = $shp.TextFrame.TextRange
$mid = $tr.Characters(1, 23)
$mid.Text = "handleAddressBookChangedEvent"

$shp = $s.Shapes.Item("TextBox 49")
$tr = $shp.TextFrame.TextRange
$mid = $tr.Characters(1, 23)
$mid.Text = "handleAddressBookChangedEvent"
